# Append one new observation record as row 36 on the "Artfynd" sheet,
# extending the used range from A1:AY35 to A1:AY36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 36

# -- Numeric columns -------------------------------------------------------
$ws.Range("A$row").Value = 111966228            # Id
$ws.Range("B$row").Value = 89183                # Taxonsorteringsordning
$ws.Range("E$row").Value = 3215                 # TaxonId
$ws.Range("Q$row").Value = 338356.4103134849    # Ost
$ws.Range("R$row").Value = 6433540.273063039    # Nord
$ws.Range("S$row").Value = 5                    # Noggrannhet

# -- Plain text columns ------------------------------------------------------
$ws.Range("C$row").Value = "Ovaliderad"                       # Valideringsstatus
$ws.Range("D$row").Value = "LC"                               # Rödlistade
$ws.Range("F$row").Value = "Rödgul trumpetsvamp"               # Artnamn
$ws.Range("G$row").Value = "Craterellus lutescens"             # Vetenskapligt namn
$ws.Range("H$row").Value = "(Fr.) Fr."                         # Auktor
$ws.Range("J$row").Value = "fruktkroppar"                      # Enhet
$ws.Range("P$row").Value = "Angertuvan, öster om, Vg"          # Lokalnamn
$ws.Range("T$row").Value = "Västra Götaland"                   # Län
$ws.Range("U$row").Value = "Ale"                               # Kommun
$ws.Range("V$row").Value = "Västergötland"                     # Provins
$ws.Range("W$row").Value = "Skepplanda"                        # Församling
$ws.Range("Z$row").Value = "00:00"                             # Starttid
$ws.Range("AB$row").Value = "00:00"                            # Sluttid
$ws.Range("AC$row").Value = "I våtmarken öster om Angertuvan. Ca 35 m söder om hyggeskanten." # Publik kommentar
$ws.Range("AH$row").Value = "Sumpskog"                         # Biotop
$ws.Range("AW$row").Value = "Thomas Grönlund"                  # Rapportör
$ws.Range("AX$row").Value = "Thomas Grönlund"                  # Observatörer

# -- Text columns that look numeric/date-like: force text storage so Excel
#    doesn't silently reinterpret them as a number or a date serial. ------
$ws.Range("I$row").NumberFormat = "@"
$ws.Range("I$row").Value = "10"                                # Antal

$ws.Range("Y$row").NumberFormat = "@"
$ws.Range("Y$row").Value = "2023-09-06"                        # Startdatum

$ws.Range("AA$row").NumberFormat = "@"
$ws.Range("AA$row").Value = "2023-09-06"                       # Slutdatum

# -- Boolean columns ---------------------------------------------------------
$ws.Range("AD$row").Value = $false   # Ej återfunnen
$ws.Range("AE$row").Value = $false   # Osäker artbestämning
$ws.Range("AG$row").Value = $false   # Ospontan

# -- Columns that exist on every row but carry no data for this record.
#    Mark them as text-formatted blanks so a cell is materialised for them,
#    mirroring the layout of the other rows in the sheet. --------------------
$ws.Range("K$row").NumberFormat = "@"    # Ålder-Stadium
$ws.Range("N$row").NumberFormat = "@"    # Metod
$ws.Range("AF$row").NumberFormat = "@"   # Bestämningsmetod
$ws.Range("AT$row").NumberFormat = "@"   # Bestämningsår
$ws.Range("AY$row").NumberFormat = "@"   # Projektnamn
